$d = $word.ActiveDocument

$d.Content.Find.Execute("N = 90,258", $true, $false, $false, $false, $false, $true, 1, $false, "N = 90,237", 2) | Out-Null
$d.Content.Find.Execute("89,943 (100)", $true, $false, $false, $false, $false, $true, 1, $false, "89,922 (100)", 2) | Out-Null
$d.Content.Find.Execute("89,612 (99)", $true, $false, $false, $false, $false, $true, 1, $false, "89,592 (99)", 2) | Out-Null
$d.Content.Find.Execute("2,017 (2.2)", $true, $false, $false, $false, $false, $true, 1, $false, "2,016 (2.2)", 2) | Out-Null
$d.Content.Find.Execute("1,284 (1.4)", $true, $false, $false, $false, $false, $true, 1, $false, "1,283 (1.4)", 2) | Out-Null
$d.Content.Find.Execute("2,680 (3.0)", $true, $false, $false, $false, $false, $true, 1, $false, "2,679 (3.0)", 2) | Out-Null
$d.Content.Find.Execute("87,578 (97)", $true, $false, $false, $false, $false, $true, 1, $false, "87,558 (97)", 2) | Out-Null
$d.Content.Find.Execute("51,733 (57)", $true, $false, $false, $false, $false, $true, 1, $false, "51,723 (57)", 2) | Out-Null
$d.Content.Find.Execute("38,525 (43)", $true, $false, $false, $false, $false, $true, 1, $false, "38,514 (43)", 2) | Out-Null
$d.Content.Find.Execute("7,120 (7.9)", $true, $false, $false, $false, $false, $true, 1, $false, "7,118 (7.9)", 2) | Out-Null
$d.Content.Find.Execute("22,117 (25)", $true, $false, $false, $false, $false, $true, 1, $false, "22,114 (25)", 2) | Out-Null
$d.Content.Find.Execute("21,321 (24)", $true, $false, $false, $false, $false, $true, 1, $false, "21,318 (24)", 2) | Out-Null
$d.Content.Find.Execute("39,700 (44)", $true, $false, $false, $false, $false, $true, 1, $false, "39,687 (44)", 2) | Out-Null
$d.Content.Find.Execute("11,586 (13)", $true, $false, $false, $false, $false, $true, 1, $false, "11,584 (13)", 2) | Out-Null
$d.Content.Find.Execute("19,663 (22)", $true, $false, $false, $false, $false, $true, 1, $false, "19,661 (22)", 2) | Out-Null
$d.Content.Find.Execute("23,668 (26)", $true, $false, $false, $false, $false, $true, 1, $false, "23,663 (26)", 2) | Out-Null
$d.Content.Find.Execute("20,901 (23)", $true, $false, $false, $false, $false, $true, 1, $false, "20,894 (23)", 2) | Out-Null
$d.Content.Find.Execute("6,074 (6.7)", $true, $false, $false, $false, $false, $true, 1, $false, "6,072 (6.7)", 2) | Out-Null
$d.Content.Find.Execute("8,366 (9.3)", $true, $false, $false, $false, $false, $true, 1, $false, "8,363 (9.3)", 2) | Out-Null
$d.Content.Find.Execute("35,546 (39)", $true, $false, $false, $false, $false, $true, 1, $false, "35,534 (39)", 2) | Out-Null
$d.Content.Find.Execute("37,147 (41)", $true, $false, $false, $false, $false, $true, 1, $false, "37,141 (41)", 2) | Out-Null
$d.Content.Find.Execute("17,052 (19)", $true, $false, $false, $false, $false, $true, 1, $false, "17,049 (19)", 2) | Out-Null
$d.Content.Find.Execute("87,401 (97)", $true, $false, $false, $false, $false, $true, 1, $false, "87,381 (97)", 2) | Out-Null
$d.Content.Find.Execute("18,206 (21)", $true, $false, $false, $false, $false, $true, 1, $false, "18,202 (21)", 2) | Out-Null
$d.Content.Find.Execute("69,526 (79)", $true, $false, $false, $false, $false, $true, 1, $false, "69,511 (79)", 2) | Out-Null
$d.Content.Find.Execute("52,105 (58)", $true, $false, $false, $false, $false, $true, 1, $false, "52,094 (58)", 2) | Out-Null
$d.Content.Find.Execute("32,075 (36)", $true, $false, $false, $false, $false, $true, 1, $false, "32,065 (36)", 2) | Out-Null
$d.Content.Find.Execute("4,961 (5.5)", $true, $false, $false, $false, $false, $true, 1, $false, "4,959 (5.5)", 2) | Out-Null
$d.Content.Find.Execute("18,247 (20)", $true, $false, $false, $false, $false, $true, 1, $false, "18,243 (20)", 2) | Out-Null
$d.Content.Find.Execute("22,703 (25)", $true, $false, $false, $false, $false, $true, 1, $false, "22,696 (25)", 2) | Out-Null
$d.Content.Find.Execute("23,634 (26)", $true, $false, $false, $false, $false, $true, 1, $false, "23,629 (26)", 2) | Out-Null
$d.Content.Find.Execute("20,713 (23)", $true, $false, $false, $false, $false, $true, 1, $false, "20,710 (23)", 2) | Out-Null
$d.Content.Find.Execute("64,281 (71)", $true, $false, $false, $false, $false, $true, 1, $false, "64,263 (71)", 2) | Out-Null
$d.Content.Find.Execute("22,781 (25)", $true, $false, $false, $false, $false, $true, 1, $false, "22,778 (25)", 2) | Out-Null
$d.Content.Find.Execute("16,297 (18)", $true, $false, $false, $false, $false, $true, 1, $false, "16,292 (18)", 2) | Out-Null
$d.Content.Find.Execute("30,604 (34)", $true, $false, $false, $false, $false, $true, 1, $false, "30,597 (34)", 2) | Out-Null
$d.Content.Find.Execute("43,357 (48)", $true, $false, $false, $false, $false, $true, 1, $false, "43,348 (48)", 2) | Out-Null
$d.Content.Find.Execute("74,546 (83)", $true, $false, $false, $false, $false, $true, 1, $false, "74,526 (83)", 2) | Out-Null
$d.Content.Find.Execute("14,991 (17)", $true, $false, $false, $false, $false, $true, 1, $false, "14,990 (17)", 2) | Out-Null
$d.Content.Find.Execute("76,549 (85)", $true, $false, $false, $false, $false, $true, 1, $false, "76,532 (85)", 2) | Out-Null
$d.Content.Find.Execute("13,285 (15)", $true, $false, $false, $false, $false, $true, 1, $false, "13,281 (15)", 2) | Out-Null
$d.Content.Find.Execute("29,461 (33)", $true, $false, $false, $false, $false, $true, 1, $false, "29,456 (33)", 2) | Out-Null
$d.Content.Find.Execute("54,052 (60)", $true, $false, $false, $false, $false, $true, 1, $false, "54,040 (60)", 2) | Out-Null
$d.Content.Find.Execute("6,745 (7.5)", $true, $false, $false, $false, $false, $true, 1, $false, "6,741 (7.5)", 2) | Out-Null
$d.Content.Find.Execute("15,370 (17)", $true, $false, $false, $false, $false, $true, 1, $false, "15,369 (17)", 2) | Out-Null
$d.Content.Find.Execute("66,477 (74)", $true, $false, $false, $false, $false, $true, 1, $false, "66,462 (74)", 2) | Out-Null
$d.Content.Find.Execute("8,411 (9.3)", $true, $false, $false, $false, $false, $true, 1, $false, "8,406 (9.3)", 2) | Out-Null

# Swap MVPA min/week Activity count <-> Machine learning labels
$d.Content.Find.Execute("MVPA min/week - Activity count, Median (Q1, Q3)", $true, $false, $false, $false, $false, $true, 1, $false, "MVPA min/week - __TEMP_SWAP__, Median (Q1, Q3)", 2) | Out-Null
$d.Content.Find.Execute("MVPA min/week - Machine learning, Median (Q1, Q3)", $true, $false, $false, $false, $false, $true, 1, $false, "MVPA min/week - Activity count, Median (Q1, Q3)", 2) | Out-Null
$d.Content.Find.Execute("MVPA min/week - __TEMP_SWAP__, Median (Q1, Q3)", $true, $false, $false, $false, $false, $true, 1, $false, "MVPA min/week - Machine learning, Median (Q1, Q3)", 2) | Out-Null
